$d = $word.ActiveDocument

# --- Paragraph 1 ---
$range = $d.Content
$found = $range.Find.Execute("במהלך התרגיל הראשון, זוהתה התקרבות של רועה צאן לשטח הגדר, ולאחר מכן נרשמה הגעה של מחבל מהשטח האדום וניסיון להברחת אמל""ח וסמים. האירוע המרכזי כלל חצייה של שני מחבלים בטיפוס על הגדר וכניסה ליישוב שקף. הכוחות פעלו במהירות וביעילות, הצליחו לעצור את האופנוען החשוד ולחסל את המחבל ששלף אקדח. כמו כן, הכוחות הצליחו לצמצם את האיום ביישוב שקף ולחסל שלושה מחבלים. עם זאת, נדרשת שיפור בתקשורת בין החמ""ל לכוחות בשטח, וביצוע סריקות מהירות יותר. כמו כן, יש להקפיד על דרישה לאמצעים רפואיים בתחילת האירוע ולא בסיומו.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found paragraph 1:" $found
if ($found) {
    $range.Text = "במהלך התרגיל הראשון, זוהתה התקרבות של רועה צאן לשטח הגדר, ולאחר מכן נרשמה הגעה של מחבל מהשטח האדום וניסיון להברחת אמל""ח וסמים. האירוע המרכזי כלל חצייה של שני מחבלים בטיפוס על הגדר וכניסה ליישוב שקף. הכוחות הגיבו במהירות, ביצעו מעצר חשוד לרועה הצאן, והצליחו לעצור את האופנוען החשוד. במהלך ההתקלות עם המחבלים, הכוחות הצליחו לחסל שלושה מחבלים, אך שני לוחמים נפצעו ושני אזרחים נהרגו. `vהכוחות פעלו היטב בזיהוי האיומים ובתגובה המהירה לאירועים. הם הצליחו לחסל את המחבלים ולמנוע נזק נוסף ליישוב. `vעם זאת, נדרשת שיפור בתקשורת בין החמ""ל לכוחות בשטח, וביצוע סריקות מהירות יותר. כמו כן, יש להקפיד על דרישה לאמצעים רפואיים בתחילת האירוע ולא בסיומו."
} else {
    Write-Host "WARNING: paragraph 1 text not found!"
}

# --- Paragraph 2 ---
$range = $d.Content
$found = $range.Find.Execute("בתרגיל השני, התקבל דיווח על חדירה בקו דיווח 455, והכוחות קפצו לנקודה ועצרו את החשוד. בהמשך, התקבל דיווח על ירי במעבר בקו דיווח 485, והכוחות צמצמו לעבר האירוע. הכוחות הצליחו לכבוש חזרה את מעבר גורביץ' ואת מוצב שקף. הכוחות פעלו בנחישות ובמהירות, הצליחו לכבוש חזרה את המוצב ולסיים את האירוע בהצלחה. עם זאת, נדרשת שיפור בתיאום בין הכוחות, במיוחד בנוגע לדיווח על מעגל פרוץ ובקשה לאמצעים רפואיים בזמן.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found paragraph 2:" $found
if ($found) {
    $range.Text = "בתרגיל השני, התקבל דיווח על חדירה בקו דיווח 455, והכוחות קפצו לנקודה ועצרו את החשוד. בהמשך, התקבל דיווח על ירי במעבר בקו דיווח 485, והכוחות צמצמו לעבר האירוע. במהלך האירוע, הכוחות הצליחו לכבוש חזרה את מעבר גורביץ' ואת מוצב שקף. `vהכוחות פעלו בצורה מקצועית בכיבוש חזרה של המעבר והמוצב, והצליחו למנוע נזק נוסף. `vעם זאת, נדרשת שיפור בתיאום בין הכוחות, במיוחד בנוגע לדיווח על מעגל פרוץ ובקשה לאמצעים רפואיים בזמן. כמו כן, יש להימנע מירי מיותר ולוודא שהמסק""ר מבצע סריקות במקומות הנכונים."
} else {
    Write-Host "WARNING: paragraph 2 text not found!"
}

# --- Paragraph 3 ---
$range = $d.Content
$found = $range.Find.Execute("התרגילים הצביעו על יכולת תגובה מהירה ונחישות של הכוחות בשטח, שהובילו להצלחות משמעותיות. עם זאת, יש לשפר את התקשורת בין החמ""ל לכוחות בשטח ואת התיאום בין הכוחות. כמו כן, יש להקפיד על בקשה לאמצעים רפואיים בתחילת האירועים ולא בסיומם. הכוחות הראו יכולת גבוהה בהתמודדות עם מצבים מורכבים, ויש להמשיך ולשפר את התהליכים הפנימיים.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found paragraph 3:" $found
if ($found) {
    $range.Text = "התרגילים הצביעו על יכולת תגובה מהירה ומקצועית של הכוחות, במיוחד בזיהוי איומים ובכיבוש חזרה של שטחים. עם זאת, יש לשפר את התקשורת בין החמ""ל לכוחות בשטח, ולהקפיד על תיאום טוב יותר בין הכוחות. הכוחות הראו יכולת גבוהה במניעת נזק נוסף, אך יש להמשיך ולשפר את התהליכים הפנימיים להבטחת הצלחה מלאה בעתיד."
} else {
    Write-Host "WARNING: paragraph 3 text not found!"
}
